$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.765.20'
$ws.Range('E2').Value = '  +0.10%  '

$ws.Range('D3').Value = '1.642.53'
$ws.Range('E3').Value = '  -0.41%  '

$ws.Range('E4').Value = '  +0.59%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '217.02'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.53%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.502'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.68%  '

$ws.Range('E7').Value = '  +0.68%  '

$ws.Range('E8').Value = '  -1.11%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.0625'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -0.74%  '

$ws.Range('E10').Value = '  -1.14%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0841'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -0.52%  '

$ws.Range('D12').Value = '1.869.14'
$ws.Range('E12').Value = '  -0.58%  '

$ws.Range('D13').Value = '1.656.67'
$ws.Range('E13').Value = '  -0.11%  '

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '4.17'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -1.32%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.526'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -1.79%  '

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '64.58'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -3.15%  '

$ws.Range('D17').Value = '26.794.44'
$ws.Range('E17').Value = '  -0.05%  '

$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  -2.60%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '214.06'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -3.35%  '

$ws.Range('E20').Value = '  +0.61%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '4.35'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -1.45%  '

$ws.Range('E22').Value = '  +12.31%  '

$ws.Range('E23').Value = '  -1.48%  '

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '9.37'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -2.49%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '144.87'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -1.70%  '

$ws.Range('E26').Value = '  +0.75%  '

$ws.Range('E27').Value = '  -2.70%  '

$ws.Range('E28').Value = '  -0.38%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '15.69'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -1.75%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.0514'
$c.Style = "Normal"

$ws.Range('E31').Value = '  +0.40%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.32'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -3.41%  '

$ws.Range('E33').Value = '  -2.77%  '

$ws.Range('D34').Value = '1.289.89'
$ws.Range('E34').Value = '  -0.40%  '

$ws.Range('E35').Value = '  -2.22%  '

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +1.21%  '

$ws.Range('E37').Value = '  -4.68%  '

$ws.Range('E38').Value = '  +2.05%  '

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.826'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.94%  '

$ws.Range('E40').Value = '  +0.63%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.809'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -0.82%  '

$ws.Range('E42').Value = '  -0.36%  '

$ws.Range('E43').Value = '  -2.16%  '

$ws.Range('D44').Value = '1.795.22'
$ws.Range('E44').Value = '  +0.17%  '

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '91.53'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -2.54%  '

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '59.89'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -0.58%  '

$ws.Range('E47').Value = '  -1.06%  '

$ws.Range('E48').Value = '  -1.82%  '

$ws.Range('E49').Value = '  +0.45%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '7.70'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -1.75%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.0976'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.64%  '

